$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add <w:noProof/> to the runs that host the three inline pictures that
#    currently lack run properties (wp14:anchorId 4AD0595B, 76831927,
#    1C99A89E -> InlineShapes #3, #4, #5 in document order).
# ---------------------------------------------------------------------------
$d.InlineShapes.Item(3).Range.NoProofing = $true
$d.InlineShapes.Item(4).Range.NoProofing = $true
$d.InlineShapes.Item(5).Range.NoProofing = $true

# ---------------------------------------------------------------------------
# 2) Split "401: 오늘의 일정이 없는 경우" into "2" + "01: 오늘의 일정이 없는
#    경우" (i.e. change the leading "4" to "2" while keeping it as its own
#    run, distinct from the rest of the sentence).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("401: 오늘의 일정이 없는 경우", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$digit = $d.Range($rng.Start, $rng.Start + 1)
$digit.Text = "2"
$digit2 = $d.Range($rng.Start, $rng.Start + 1)
$digit2.Font.Bold = $true
$digit2.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3) Merge the three runs "을 방지하기 위해서 " + "100" + "자 이상 적은 메모에
#    대해서만 입력 데이터를 생성함. 따라서 이를 사전에 충분히 고지해야 함."
#    into a single run of text.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("을 방지하기 위해서 100자 이상 적은 메모에 대해서만 입력 데이터를 생성함. 따라서 이를 사전에 충분히 고지해야 함.", $false, $false, $false, $false, $false, $true, 1, $false, "을 방지하기 위해서 100자 이상 적은 메모에 대해서만 입력 데이터를 생성함. 따라서 이를 사전에 충분히 고지해야 함.", 2)
